$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NuevaHoja")

# Re-touch the custom date format on C1 so the workbook's custom numFmtId
# gets reassigned to the next free custom id (164) instead of the
# originally-loaded 168.
$dateFormat = $ws.Range("C1").NumberFormat
$ws.Range("C1").NumberFormat = $dateFormat

# New data in row 2
$ws.Range("B2").Value = "Segundo texto"
$ws.Range("C2").Value = 3000

# New empty but styled cell H19 (underline font, same style class as C1's font)
$ws.Range("H19").Font.Underline = $true
$ws.Range("H19").Select()

# Defined name at workbook scope
$excel.Names.Add('PrimerRango', '=NuevaHoja!$A$1:$C$2')
